# feat: add 2022-Q4 data
#
# 1) Insert a new "2022-Q4" worksheet (cloned from "2022-Q3" so it picks up
#    the same layout/styles) positioned right after "总计" and before the
#    existing "2022-Q3" sheet, then fill it with the new quarter's fund data.
# 2) Insert a new top data row in "总计" for 2022-Q4 and renumber the
#    existing rows' index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Part 1: new "2022-Q4" sheet
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$zj = $wb.Worksheets.Item("总计")
$q3.Copy($null, $zj)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Row 2: 001628 招商体育文化休闲股票A
$q4.Cells.Item(2,2).Value = "'001628"
$q4.Cells.Item(2,3).Value = "'招商体育文化休闲股票A"
$q4.Cells.Item(2,4).Value = "'2.33"
$q4.Cells.Item(2,5).Value = "'93.03"
$q4.Cells.Item(2,6).Value = "'5.41"
$q4.Cells.Item(2,7).Value = "'0.1261"
$q4.Range("B2:G2").ClearFormats()
$q4.Cells.Item(2,8).Value = 3

# Row 3: 015395 招商体育文化休闲股票C
$q4.Cells.Item(3,2).Value = "'015395"
$q4.Cells.Item(3,3).Value = "'招商体育文化休闲股票C"
$q4.Cells.Item(3,4).Value = "'0.29"
$q4.Cells.Item(3,5).Value = "'93.03"
$q4.Cells.Item(3,6).Value = "'5.41"
$q4.Cells.Item(3,7).Value = "'0.0157"
$q4.Range("B3:G3").ClearFormats()
$q4.Cells.Item(3,8).Value = 3

# Row 4 (brand new row): 165531 信诚多策略灵活配置混合（LOF）
# Copy column-A formatting from the row above so the new index cell gets the
# same bordered/bold style as the rest of the column.
$q4.Cells.Item(3,1).Copy()
$q4.Cells.Item(4,1).PasteSpecial(-4122)
$q4.Cells.Item(4,1).Value = 2
$q4.Cells.Item(4,2).Value = "'165531"
$q4.Cells.Item(4,3).Value = "'信诚多策略灵活配置混合（LOF）"
$q4.Cells.Item(4,4).Value = "'0.89"
$q4.Cells.Item(4,5).Value = "'72.25"
$q4.Cells.Item(4,6).Value = "'1.01"
$q4.Cells.Item(4,7).Value = "'0.0090"
$q4.Range("B4:G4").ClearFormats()
$q4.Cells.Item(4,8).Value = 10

# ---------------------------------------------------------------------
# Part 2: update the "总计" summary sheet with the new 2022-Q4 row
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("总计")
$ws.Rows.Item(2).Insert()

$ws.Cells.Item(2,2).Value = "'2022-Q4"
$ws.Cells.Item(2,3).Value = 3
$ws.Cells.Item(2,4).Value = 0.15
$ws.Range("B2:D2").ClearFormats()

# New row-2 index cell: copy formatting from the (now shifted) row below.
$ws.Cells.Item(3,1).Copy()
$ws.Cells.Item(2,1).PasteSpecial(-4122)
$ws.Cells.Item(2,1).Value = 0

# Renumber the index column for the rows that shifted down one position.
$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(5,1).Value = 3
